$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (id=1)
$ws.Range("B2").Value = 4.511983941394703
$ws.Range("D2").Value = 7.999999999999998
$ws.Range("F2").Value = 4.999999999999999
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 2.688252163597601
$ws.Range("L2").Value = 2.359683601981216

# Row 3 (id=2)
$ws.Range("B3").Value = 3.805838572905071
$ws.Range("D3").Value = 3.223773863882188
$ws.Range("F3").Value = 4.387903281927955

# Row 4 (id=3)
$ws.Range("B4").Value = 5.050287379341913
$ws.Range("D4").Value = 6.084029194605166
$ws.Range("F4").Value = 3.592697288672839
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 5.474135654747734
